$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 119, shifting existing rows 119..207 down to 120..208
$ws.Rows.Item(119).Insert()

# Populate the new row 119 with the new weekly price entry
$ws.Range("A119").Value = 4
$ws.Range("B119").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C119").Value = "Los Lagos"
$ws.Range("D119").Value = 44574
$ws.Range("E119").Value = 10
$ws.Range("F119").Value = 100112043
$ws.Range("G119").Value = "Pepino ensalada"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 200
$ws.Range("K119").Value = 14000
$ws.Range("L119").Value = 14000
$ws.Range("M119").Value = 14000
$ws.Range("N119").Value = "$/caja 60 unidades"
$ws.Range("O119").Value = "Región de Arica y Parinacota"
$ws.Range("P119").Value = 233
$ws.Range("Q119").Value = 60
$ws.Range("R119").Value = "Hortaliza"
